$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cases tab query (B2): append an ORDER BY / LIMIT clause.
$caseQuery = $ws.Range("B2").Value()
$ws.Range("B2").Value = $caseQuery + "`n order By ss.study_subject_id ASC LIMIT 100 "

# Samples tab query (B3): append an ORDER BY / LIMIT clause.
$sampleQuery = $ws.Range("B3").Value()
$ws.Range("B3").Value = $sampleQuery + "`n order By samp.sample_id ASC LIMIT 100"

# Files tab query (B4): append an ORDER BY / LIMIT clause.
$fileQuery = $ws.Range("B4").Value()
$ws.Range("B4").Value = $fileQuery + "`n order By f.file_name ASC LIMIT 100"

# Selection moved from C4 to B4 (and the view's frozen/scrolled top-left cell resets).
$ws.Range("B4").Select()
